$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Humans or Sentient beings"
$ws.Range("B2").Value = 0.453581378758874
$ws.Range("K2").Value = 0.295505791021036
$ws.Range("L2").Value = 0.498063798608243
$ws.Range("N2").Value = 0.434044780311437

# Row 3 - "Fellow citizens"
$ws.Range("B3").Value = 0.318130319930611
$ws.Range("K3").Value = 0.451904268645967
$ws.Range("L3").Value = 0.227706886529004
$ws.Range("N3").Value = 0.338155063004691

# Row 4 - "Family and self"
$ws.Range("B4").Value = 0.172200105515648
$ws.Range("E4").Value = 0.16615846782681
$ws.Range("K4").Value = 0.19166544338377
$ws.Range("L4").Value = 0.194374624175739
$ws.Range("N4").Value = 0.175591451815085
